$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Copy the header formatting (bold, border, centered) from the Weekly sheet's
# header row onto the new header row, and the date formatting from its date
# column onto the new date column so the same styles get reused.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows ---
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 960
$wsForecast.Range("C2").Value = 959.9982249991897
$wsForecast.Range("D2").Value = 959.9982250007456

$wsForecast.Range("A3").Value = 44948.99999999999
$wsForecast.Range("B3").Value = 20
$wsForecast.Range("C3").Value = 19.99790499555152
$wsForecast.Range("D3").Value = 19.9979049971259

$wsForecast.Range("A4").Value = 44955.99999999999
$wsForecast.Range("B4").Value = 0
$wsForecast.Range("C4").Value = -450.0022573068455
$wsForecast.Range("D4").Value = -450.0022526964157

$wsForecast.Range("A5").Value = 44962.99999999999
$wsForecast.Range("B5").Value = 0
$wsForecast.Range("C5").Value = -920.0024231122477
$wsForecast.Range("D5").Value = -920.0024074397342

$wsForecast.Range("A6").Value = 44969.99999999999
$wsForecast.Range("B6").Value = 0
$wsForecast.Range("C6").Value = -1390.00259222213
$wsForecast.Range("D6").Value = -1390.002559551571

$wsForecast.Range("A7").Value = 44976.99999999999
$wsForecast.Range("B7").Value = 0
$wsForecast.Range("C7").Value = -1860.002762306502
$wsForecast.Range("D7").Value = -1860.002709830686

$wsForecast.Range("A8").Value = 44983.99999999999
$wsForecast.Range("B8").Value = 0
$wsForecast.Range("C8").Value = -2330.002932619455
$wsForecast.Range("D8").Value = -2330.002859674468

$wsForecast.Range("A9").Value = 44990.99999999999
$wsForecast.Range("B9").Value = 0
$wsForecast.Range("C9").Value = -2800.00310565673
$wsForecast.Range("D9").Value = -2800.003008752886

$wsForecast.Range("A10").Value = 44997.99999999999
$wsForecast.Range("B10").Value = 0
$wsForecast.Range("C10").Value = -3270.003277581429
$wsForecast.Range("D10").Value = -3270.003155023861

$wsForecast.Range("A11").Value = 45004.99999999999
$wsForecast.Range("B11").Value = 0
$wsForecast.Range("C11").Value = -3740.003454056607
$wsForecast.Range("D11").Value = -3740.00329962216

$wsForecast.Range("A1").Select() | Out-Null

# Restore the originally active sheet/tab so the workbook view state
# (activeTab) is unaffected by having created and populated the new sheet.
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
